$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the list date in D1 (new KB-Liste date)
$ws.Range("D1").Value = 44407

# 2. Row 32: "Sunboy-ET" becomes "Sunboy-ET (Prüf)" with " (Prüf)" bold,
#    and the cell takes on the shaded/bold-font style used by other
#    "(Prüf)" rows (copy formatting from A30, which already carries it).
$ws.Range("A30").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A32").Value = "Sunboy-ET (Prüf)"
$ws.Range("A32").Characters(10, 7).Font.Bold = $true

# 3. Row 38 (X-Amnesty-ET): clear the colour/marking cell in column D
$ws.Range("D38").ClearContents()

# 4. Row 39 (X-Andero SG-ET): annotate the remark with "(nicht 004/005)" in bold
$ws.Range("E39").Value = "Inhalt, Euter, Nutzungsdau (nicht 004/005)"
$ws.Range("E39").Characters(28, 15).Font.Bold = $true

# 5. Row 41 (X-Arland SG-ET): annotate the remark with "(nicht 011/012)" in bold
$ws.Range("E41").Value = " +1627, Inhalt, Nutzungsdau (nicht 011/012)"
$ws.Range("E41").Characters(29, 15).Font.Bold = $true

# 6. Insert a new row for "X-SUV PP" before the X-Swingman-ET row (old row 46)
$ws.Rows(46).Insert()
$ws.Range("A45:F45").Copy()
$ws.Range("A46:F46").PasteSpecial(-4122)
$ws.Range("A46").Value = "X-SUV PP"
$ws.Range("B46").Value = "RH Op"
$ws.Range("C46").Value = "67 o"
$ws.Range("D46").Value = "rot/dünn"
$ws.Range("E46").Value = "Roboter, Geburt, ZZ, Fruchtbar, Exterieur"
$ws.Range("F46").Value = " -"

# 7. Move the active-cell selection to D2 (matches the refreshed sheetView)
$ws.Range("D2").Select()
